$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two additional dates (columns F-I)
$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

# Copy the header style/formatting (bold, border, centered) from an existing
# header cell onto the new header cells so they share the same cell style.
$ws.Range("D1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate default attendance data ("A" / "00:00:00") for every student row
# in the two newly added date columns.
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 6).Value = "A"
    $ws.Cells.Item($r, 7).Value = "00:00:00"
    $ws.Cells.Item($r, 8).Value = "A"
    $ws.Cells.Item($r, 9).Value = "00:00:00"
}

# Row 27 (ABEESHA ESHAL) was marked present on 07-04-2025.
$ws.Cells.Item(27, 3).Value = 1
$ws.Cells.Item(27, 4).Value = "P"
$ws.Cells.Item(27, 5).Value = "09:46:00 AM"
